$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = "MODELO3"
$ws.Range("E14").Value = "MODELO2"

# Row 15
$ws.Range("A15").Value = "base"
$ws.Range("B15").Value = 154
$ws.Range("C15").Value = "poliacetal"
$ws.Range("E15").Value = 133

# Row 16
$ws.Range("A16").Value = "lateral"
$ws.Range("B16").Value = 104
$ws.Range("C16").Value = "poliacetal"
$ws.Range("E16").Value = 71
$ws.Range("F16").Value = 135
$ws.Range("G16").Value = "(al)"

# Row 17
$ws.Range("A17").Value = "arma"
$ws.Range("B17").Value = 117
$ws.Range("C17").Value = "aço"
$ws.Range("E17").Value = 271

# Row 18
$ws.Range("A18").Value = "eixo "
$ws.Range("B18").Value = 20
$ws.Range("C18").Value = "aço"
$ws.Range("E18").Value = 39

# Row 19
$ws.Range("A19").Value = "polia"
$ws.Range("B19").Value = 25
$ws.Range("C19").Value = "poliacetal"
$ws.Range("D19").Formula = "=SUM(B13:B24)"
$ws.Range("E19").Value = 25

# Row 20
$ws.Range("A20").Value = "rolamento"
$ws.Range("B20").Value = 15
$ws.Range("C20").Value = "aço"
$ws.Range("E20").Value = 30

# Row 21
$ws.Range("A21").Value = "tampa"
$ws.Range("B21").Value = 154
$ws.Range("C21").Value = "poliacetal"
$ws.Range("E21").Value = 144

# Row 22
$ws.Range("A22").Value = "mancal"
$ws.Range("E22").Value = 18
$ws.Range("F22").Value = 34
$ws.Range("G22").Value = 97

# Row 23
$ws.Range("A23").Value = "mancal"
$ws.Range("E23").Value = 18
$ws.Range("F23").Value = 34
$ws.Range("G23").Value = 97

# Row 26
$ws.Range("E26").Formula = "=SUM(B13,E15:E23)"

# Window/view state
$ws.Range("C24").Select() | Out-Null
